$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Mission numbers for rows 62-66 (missions 60-64)
$ws1.Cells.Item(62, 1).Value = 60
$ws1.Cells.Item(63, 1).Value = 61
$ws1.Cells.Item(64, 1).Value = 62
$ws1.Cells.Item(65, 1).Value = 63
$ws1.Cells.Item(66, 1).Value = 64

# Row 62 - mission 60 (N=14 microSWIFTs Deployed, O=15 microSWIFTs Retrieved, Q=17 Start Time, R=18 End Time)
$ws1.Cells.Item(62, 14).Value = "56,57,60,36,30,45,71,58"
$ws1.Cells.Item(62, 15).Value = "56,57,60,36,30,45,71,58"
$ws1.Cells.Item(62, 17).Value = "2021-10-27T14:30:00"
$ws1.Cells.Item(62, 18).Value = "2021-10-27T14:50:00"

# Row 63 - mission 61
$ws1.Cells.Item(63, 14).Value = "31,33,34,35,37,38,40,41"
$ws1.Cells.Item(63, 15).Value = "31,33,34,35,37,38,40,41"
$ws1.Cells.Item(63, 17).Value = "2021-10-27T15:16:00"
$ws1.Cells.Item(63, 18).Value = "2021-10-27T15:20:00"

# Row 64 - mission 62
$ws1.Cells.Item(64, 14).Value = "3,12,13,16,17,18,19,20,21,23,24,26,27,29,31,33,34,35,37,38,40,41,42,43,48,49,50,59,72,73"
$ws1.Cells.Item(64, 15).Value = "3,12,13,16,17,18,19,20,21,23,24,26,27,29,31,33,34,35,37,38,40,41,42,43,48,49,50,59,72,73"
$ws1.Cells.Item(64, 17).Value = "2021-10-27T15:32:00"
$ws1.Cells.Item(64, 18).Value = "2021-10-27T15:48:00"

# Row 65 - mission 63
$ws1.Cells.Item(65, 14).Value = "3,71,73,56,72,59,60,10,58,57,12,13,16,17,18,20,21,23,24,26,27,29,30,31,33,34,35,36,37,38,40,41,42,43,45,48,49,50"
$ws1.Cells.Item(65, 15).Value = "3,71,73,56,72,59,60,10,58,57,12,13,16,17,18,20,21,23,24,26,27,29,30,31,33,34,35,36,37,38,40,41,42,43,45,48,49,50"
$ws1.Cells.Item(65, 17).Value = "2021-10-27T16:10:00"
$ws1.Cells.Item(65, 18).Value = "2021-10-27T16:35:00"

# Row 66 - mission 64
$ws1.Cells.Item(66, 14).Value = "3,56,57,58,59,60,71,72,73,10,12,13,16,17,18,20,21,23,24,26,27,29,30,31,33,34,35,36,37,38,40,41,42,43,45,48,49,50"
$ws1.Cells.Item(66, 15).Value = "3,56,57,58,59,60,71,72,73,10,12,13,16,17,18,20,21,23,24,26,27,29,30,31,33,34,35,36,37,38,40,41,42,43,45,48,49,50"
$ws1.Cells.Item(66, 17).Value = "2021-10-27T16:29:00"
$ws1.Cells.Item(66, 18).Value = "2021-10-27T16:39:00"

# Row heights grow because columns N/O wrap the (longer) microSWIFT lists.
$ws1.Rows.Item(62).RowHeight = 17
$ws1.Rows.Item(63).RowHeight = 17
$ws1.Rows.Item(64).RowHeight = 68
$ws1.Rows.Item(65).RowHeight = 68
$ws1.Rows.Item(66).RowHeight = 68

# Update sheet views / selections to match the final state recorded in the file:
# Sheet2 selection moves to H7 and it is no longer the active tab.
$ws2.Select()
$ws2.Range("H7").Select()

# Sheet1 becomes the active tab with cell S66 selected (last edited cell).
$ws1.Select()
$ws1.Range("S66").Select()
